# Add the 27th (4th) class date/attendance column to the "November" sheet,
# then recalc and leave the "Overall Attendance" sheet active (as the author did
# before saving).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("November")

# Header: new date (27) and its time slot, mirroring columns D/E/F.
$ws.Range("G10").Value = 27
$ws.Range("G11").Value = "09:00 TO`n12:00"

# "Total Classes" row - 3 classes were held on this date.
$ws.Range("G12").Value = 3

# Per-student attendance for the new date (mirrors D/E/F pattern: 3 = attended,
# 0 = absent). Row 17 (S# 4) missed this date; rows 21 and 26 were already
# absent for the whole month.
$ws.Range("G14").Value = 3
$ws.Range("G15").Value = 3
$ws.Range("G16").Value = 3
$ws.Range("G17").Value = 0
$ws.Range("G18").Value = 3
$ws.Range("G19").Value = 3
$ws.Range("G20").Value = 3
$ws.Range("G21").Value = 0
$ws.Range("G22").Value = 3
$ws.Range("G23").Value = 3
$ws.Range("G24").Value = 3
$ws.Range("G25").Value = 3
$ws.Range("G26").Value = 0

$excel.Calculate()

# Mirror the author's final UI state: last interaction was on November (cell
# G24), then they switched to "Overall Attendance" before saving.
$ws.Range("G24").Select()
$overall = $wb.Worksheets.Item("Overall Attendance")
$overall.Activate()
$overall.Range("G11").Select()
